# A new September transaction ("balance your axis", logged 2024-09-10
# 13:32:42) was recorded at the top of the "2024" worksheet's September
# list. Inserting a new row at row 36 pushes the existing September list
# (R:S), the trailing August list (P:Q) and the "Broadband" label (A)
# down by one row each - exactly matching how this running log grows.
# The sheet's dimension grows from A1:Y127 to A1:Y128 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows("36:36").Insert()

$ws.Range("R36").Value = "balance your axis"
$ws.Range("S36").Value = "2024-09-10 13:32:42"
